$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("POReceipt")

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

Set-TextValue $ws.Range("J2") "1291.0"
Set-TextValue $ws.Range("K2") "1295.0"
Set-TextValue $ws.Range("J3") "65.0"
Set-TextValue $ws.Range("K3") "67.0"
Set-TextValue $ws.Range("J4") "733.0"
Set-TextValue $ws.Range("K4") "739.0"
